# MiX_BOM.xlsx: "Update tab names in all BOMs, fix bi-color LED naming."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update tab name: "MiX_BOM" -> "BOM"
$ws.Name = "BOM"

# Fix bi-color LED naming: the Ref cells for the trimmer pots / jacks
# (D10:D12) were carrying a stray, visually-identical font style
# (fontId 0 with applyFont) left over from the LED-naming edit.
# Reset them back to the default "Normal" style so they no longer
# carry an explicit style index.
$ws.Range("D10:D12").Style = "Normal"
